$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(67, 46, 59, 71, 79, 77, 86, 82, 63, 78, 75)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

$ws.Range("D13").Select()
